$d = $word.ActiveDocument

function Find-ParaByText($needle) {
    foreach ($p in $d.Paragraphs) {
        if ($p.Range.Text -like $needle) {
            return $p
        }
    }
    return $null
}

# Locate the part3 header paragraph ({{part3.title}}) and the paragraph
# right before it (the page break that separates part2 from part3). New
# content is inserted after that page-break paragraph, whose run carries no
# bold formatting, so the freshly split-off paragraphs stay unformatted.
$titlePara = Find-ParaByText("*{{part3.title}}*")
$anchor = $titlePara.Previous()

$newLines = @(
    "{{FOR s3 IN part3.senderLines}}",
    "{{INS `$s3}}",
    "{{END-FOR s3}}",
    "",
    "{{FOR a3 IN part3.addresseeLines}}",
    "{{INS `$a3}}",
    "{{END-FOR a3}}",
    "Datum: {{part3.dateLine}}",
    ""
)

foreach ($line in $newLines) {
    $anchor.Range.InsertParagraphAfter() | Out-Null
    $anchor = $anchor.Next()
    $anchor.Range.Text = $line
}

# Give the date line its characteristic right-aligned-block indent, matching
# the "Datum:" lines used elsewhere in this template. Re-locate it fresh
# (by its distinctive placeholder) now that the document has shifted.
$datePara = Find-ParaByText("*part3.dateLine*")
$datePara.Range.ParagraphFormat.LeftIndent = 319.05

# Refine the bold header line itself: {{part3.title}} -> Betreff: {{part3.subject}}
# Re-locate it fresh since earlier inserts shifted paragraph ranges around.
$titlePara = Find-ParaByText("*{{part3.title}}*")
$titlePara.Range.Find.Execute("{{part3.title}}", $true, $false, $false, $false,
                               $false, $true, 1, $false,
                               "Betreff: {{part3.subject}}", 2) | Out-Null
